# Automatic update of files.
#
# Inserts a new 'Knärot – ekologi samt krav på livsmiljön' section (heading,
# five body paragraphs, a 'Referenser - knärot' heading and six reference
# paragraphs) directly before the existing 'Järpe – ekologi samt krav på
# livsmiljön' heading, and bumps the date stamped in the document's first-page
# header from 2023-09-13 to 2023-09-15.

$d = $word.ActiveDocument

# Locate the anchor paragraph: the existing Heading 1 paragraph that starts
# the "Järpe – ekologi..." section. The new "Knärot" section is inserted
# directly before it.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Järpe – ekologi") -and $p.Style.NameLocal -eq "Heading 1") {
        $anchor = $p
        break
    }
}
if ($anchor -eq $null) { throw "Anchor paragraph (Järpe heading) not found" }

# Insert the new paragraphs in reverse order: each InsertParagraphBefore()
# call turns $anchor into the freshly-created blank paragraph positioned
# right before whatever $anchor last pointed to, so filling it in and
# repeating while walking the source paragraphs back-to-front reconstructs
# the forward (top-to-bottom) order in the final document.
#
# Italic is explicitly set (0 or 1) on the new paragraph and on every run
# typed into it, because Word's COM layer otherwise lets ambient 'current
# typing formatting' (e.g. italic left over from a previous InsertAfter)
# leak into text inserted later in the script.

# paragraph 13/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'SLU Artdatabanken, 2021. '
$rLen = 25
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$rLen = 36
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('SLU Artdatabanken, Uppsala ')
$rLen = 27
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 12/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Skogsstyrelsen, 2022. '
$rLen = 22
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Vägledning för hänsyn till knärot. ')
$rLen = 35
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$rLen = 128
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 11/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. '
$rLen = 54
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$rLen = 67
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$rLen = 38
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 10/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. '
$rLen = 117
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$rLen = 90
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$rLen = 36
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 9/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. '
$rLen = 62
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$rLen = 114
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$rLen = 39
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 8/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'de Graaf M & Roberts M.R., 2009. '
$rLen = 33
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$rLen = 80
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$rLen = 44
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 7/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Heading 2')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Referenser - knärot'
$rLen = 19
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 6/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$rLen = 868
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 5/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$rLen = 1337
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 4/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: '
$rLen = 205
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$rLen = 865
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen

# paragraph 3/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Samuel Johnsons doktorsavhandling '
$rLen = 34
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$rLen = 82
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$rLen = 162
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$rLen = 205
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('Vidare ')
$rLen = 7
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen
$ip = $d.Range($cursor, $cursor)
$ip.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$rLen = 118
$d.Range($cursor, $cursor + $rLen).Font.Italic = 1
$cursor = $cursor + $rLen

# paragraph 2/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Normal')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$rLen = 684
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# paragraph 1/13 (final top-to-bottom order)
$anchor.Range.InsertParagraphBefore()
$anchor.Font.Italic = 0
$anchor.Style = $d.Styles.Item('Heading 1')
$pStart = $anchor.Range.Start
$cursor = $pStart
$anchor.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'
$rLen = 40
$d.Range($cursor, $cursor + $rLen).Font.Italic = 0
$cursor = $cursor + $rLen

# --- Update the date stamped in the first-page header -------------------
# (rsids differ per part, so find the header whose text contains the date
# rather than assuming a fixed Headers collection index).
$sec = $d.Sections.Item(1)
$headers = $sec.Headers
$dateHeader = $null
for ($i = 1; $i -le $headers.Count; $i++) {
    $h = $headers.Item($i)
    if ($h.Exists -and $h.Range.Text.Contains("2023-09-13")) {
        $dateHeader = $h
        break
    }
}
if ($dateHeader -eq $null) { throw "Header containing the date was not found" }
$dateHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Host "Inserted Knärot section and updated header date."
